$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19792.54129835656
$ws.Range("C2").Value = 259490.0913668228
$ws.Range("D2").Value = 7.627474788768424

$ws.Range("B3").Value = 15066.47555447738
$ws.Range("C3").Value = 140979.0787425693
$ws.Range("D3").Value = 10.68702937262704

$ws.Range("B4").Value = 621.6720529563064
$ws.Range("C4").Value = 6179.869957071091
$ws.Range("D4").Value = 10.05963001284486

$ws.Range("B5").Value = 30290.74036012588
$ws.Range("C5").Value = 316598.7517522383
$ws.Range("D5").Value = 9.567548890347679

$ws.Range("B6").Value = 5701.982570874461
$ws.Range("C6").Value = 117665.1833764828
$ws.Range("D6").Value = 4.845938626237748
